{"js": "// Add two new paragraphs right after the final KNN illustration image and\n// before the trailing blank paragraph, describing how to pick k and when\n// KNN is efficient to use.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate every paragraph's inline pictures so we can find the paragraph\n// that holds the last image in the document (the KNN diagram).\nfor (const paragraph of paragraphs.items) {\n  paragraph.inlinePictures.load(\"items\");\n}\nawait context.sync();\n\nlet lastImageParagraphIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].inlinePictures.items.length > 0) {\n    lastImageParagraphIndex = i;\n  }\n}\n\nif (lastImageParagraphIndex === -1) {\n  throw new Error(\"Could not find the KNN image paragraph to anchor the new text on.\");\n}\n\n// The new paragraphs are inserted immediately after the image paragraph,\n// i.e. immediately before the paragraph that currently follows it.\nconst anchorParagraph = paragraphs.items[lastImageParagraphIndex + 1];\n\nconst newParagraphTexts = [\n  \"While selecting value of k, it can be looked at the square root of n where n is the total number of data points. Make sure that k is odd to avoid confusion between two classes of data. \",\n  \"KNN algorithm is efficient to use when the dataset is small, data is noise free (no random data in specific field) and it is labeled. \"\n];\n\nfor (const text of newParagraphTexts) {\n  const inserted = anchorParagraph.insertParagraph(text, Word.InsertLocation.before);\n  inserted.alignment = Word.Alignment.justified;\n}\n\nawait context.sync();\n", "ps1": "# Add two new paragraphs right after the final KNN illustration image and\n# before the trailing blank paragraph, describing how to pick k and when\n# KNN is efficient to use.\n$d = $word.ActiveDocument\n\n# Find the last inline image in the document (the KNN diagram) and figure\n# out which paragraph contains it.\n$shapeCount = $d.InlineShapes.Count\nif ($shapeCount -eq 0) {\n    throw \"No inline images found to anchor the new text on.\"\n}\n$lastShape = $d.InlineShapes.Item($shapeCount)\n$shapeStart = $lastShape.Range.Start\n\n$imgParaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($shapeStart -ge $p.Range.Start -and $shapeStart -lt $p.Range.End) {\n        $imgParaIndex = $i\n        break\n    }\n}\nif ($imgParaIndex -eq -1) {\n    throw \"Could not locate the paragraph containing the KNN image.\"\n}\n\n# The new paragraphs go immediately after the image paragraph, i.e.\n# immediately before the paragraph that currently follows it.\n$anchorIndex = $imgParaIndex + 1\n$anchor = $d.Paragraphs.Item($anchorIndex)\n$rng = $anchor.Range\n$rng.Collapse(1)  # wdCollapseStart\n$rng.InsertParagraphBefore()\n$rng.InsertParagraphBefore()\n\n$para1 = $d.Paragraphs.Item($anchorIndex)\n$para2 = $d.Paragraphs.Item($anchorIndex + 1)\n$para1.Range.Text = \"While selecting value of k, it can be looked at the square root of n where n is the total number of data points. Make sure that k is odd to avoid confusion between two classes of data. \"\n$para2.Range.Text = \"KNN algorithm is efficient to use when the dataset is small, data is noise free (no random data in specific field) and it is labeled. \"\n\n$d.Save()\n"}
